$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.000.04'
$ws.Range("E2").Value = '  +2.96%  '
$ws.Range("D3").Value = '3.599.93'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '657.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.70'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +16.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.422'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.37%  '
$ws.Range("E9").Value = '  +6.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").Value = '3.594.55'
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.87%  '
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").Value = '97.923.86'
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("D16").Value = '4.269.44'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("E17").Value = '  +3.45%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.591.26'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.521'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +10.81%  '
$ws.Range("E23").Value = '  +2.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '516.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000205'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.19%  '
$ws.Range("D29").Value = '3.794.12'
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.159'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +13.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.185'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.89%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '617.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.12%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.568'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.98'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +14.19%  '
$ws.Range("E42").Value = '  +3.02%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.921'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0440'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.31'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.78%  '
$ws.Range("B50").Value = 'MantraDAO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.13%  '
